# Auto-generated Excel COM-interop edit script
# Refreshes per-cell market-price figures (columns H-N) across the
# Odin_Profits leve-profit sheets, matching the scheduled-runner update.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 693.1818
$ws.Range("J19").Value = 388.33334
$ws.Range("L19").Value = 388.33334
$ws.Range("N19").Value = -738.33334
$ws.Range("H88").Value = 2972.9167
$ws.Range("I88").Value = 594
$ws.Range("J88").Value = 3189.182
$ws.Range("K88").Value = 594
$ws.Range("L88").Value = 3189.182
$ws.Range("M88").Value = -188
$ws.Range("N88").Value = -4001.182
$ws.Range("H91").Value = 2972.9167
$ws.Range("I91").Value = 594
$ws.Range("J91").Value = 3189.182
$ws.Range("K91").Value = 594
$ws.Range("L91").Value = 3189.182
$ws.Range("M91").Value = 810
$ws.Range("N91").Value = -5997.182
$ws.Range("H112").Value = 3025.1428
$ws.Range("I112").Value = 4030.6667
$ws.Range("J112").Value = 2750.9092
$ws.Range("K112").Value = 12092.0001
$ws.Range("L112").Value = 8252.7276
$ws.Range("M112").Value = -10984.0001
$ws.Range("N112").Value = -10468.7276
$ws.Range("H137").Value = 3899.742
$ws.Range("J137").Value = 3350.5
$ws.Range("L137").Value = 10051.5
$ws.Range("N137").Value = -15151.5
$ws.Range("H138").Value = 3660.2
$ws.Range("I138").Value = 2130.8
$ws.Range("J138").Value = 4233.725
$ws.Range("K138").Value = 6392.400000000001
$ws.Range("L138").Value = 12701.175
$ws.Range("M138").Value = -1252.400000000001
$ws.Range("N138").Value = -22981.175

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2784611.8
$ws.Range("I32").Value = 6273.6855
$ws.Range("K32").Value = 6273.6855
$ws.Range("M32").Value = -5986.6855
$ws.Range("H61").Value = 4231.6855
$ws.Range("I61").Value = 3760.439
$ws.Range("K61").Value = 3760.439
$ws.Range("M61").Value = -3548.439
$ws.Range("H76").Value = 35910.625
$ws.Range("J76").Value = 35910.625
$ws.Range("L76").Value = 35910.625
$ws.Range("N76").Value = -36586.625
$ws.Range("H79").Value = 35910.625
$ws.Range("J79").Value = 35910.625
$ws.Range("L79").Value = 35910.625
$ws.Range("N79").Value = -38250.625
$ws.Range("H125").Value = 42333
$ws.Range("J125").Value = 42333
$ws.Range("L125").Value = 42333
$ws.Range("N125").Value = -52173
$ws.Range("H136").Value = 4231.6855
$ws.Range("I136").Value = 3760.439
$ws.Range("K136").Value = 11281.317
$ws.Range("M136").Value = -8731.316999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 7539.385
$ws.Range("I99").Value = 7055.304
$ws.Range("J99").Value = 8711.368
$ws.Range("K99").Value = 7055.304
$ws.Range("L99").Value = 8711.368
$ws.Range("M99").Value = -5557.304
$ws.Range("N99").Value = -11707.368

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 216.83333
$ws.Range("I7").Value = 160.2
$ws.Range("K7").Value = 160.2
$ws.Range("M7").Value = -47.19999999999999
$ws.Range("H31").Value = 10424301
$ws.Range("I31").Value = 35728956
$ws.Range("J31").Value = 4737.7646
$ws.Range("K31").Value = 35728956
$ws.Range("L31").Value = 4737.7646
$ws.Range("M31").Value = -35728661
$ws.Range("N31").Value = -5327.7646
$ws.Range("H34").Value = 10424301
$ws.Range("I34").Value = 35728956
$ws.Range("J34").Value = 4737.7646
$ws.Range("K34").Value = 35728956
$ws.Range("L34").Value = 4737.7646
$ws.Range("M34").Value = -35728754
$ws.Range("N34").Value = -5141.7646
$ws.Range("H107").Value = 294.1579
$ws.Range("I107").Value = 233.94118
$ws.Range("J107").Value = 806
$ws.Range("K107").Value = 233.94118
$ws.Range("L107").Value = 806
$ws.Range("M107").Value = 1686.05882
$ws.Range("N107").Value = -4646
$ws.Range("H132").Value = 11552.333
$ws.Range("I132").Value = 5929.3335
$ws.Range("K132").Value = 17788.0005
$ws.Range("M132").Value = -15258.0005
$ws.Range("H134").Value = 66679440
$ws.Range("I134").Value = 111117760
$ws.Range("K134").Value = 333353280
$ws.Range("M134").Value = -333350745

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1387.5555
$ws.Range("J34").Value = 3736.6667
$ws.Range("L34").Value = 11210.0001
$ws.Range("N34").Value = -11378.0001
$ws.Range("H37").Value = 111914.27
$ws.Range("J37").Value = 111914.27
$ws.Range("L37").Value = 335742.81
$ws.Range("N37").Value = -335966.81
$ws.Range("H52").Value = 2319022.8
$ws.Range("J52").Value = 2319022.8
$ws.Range("L52").Value = 6957068.399999999
$ws.Range("N52").Value = -6957600.399999999
$ws.Range("H113").Value = 3290440.8
$ws.Range("I113").Value = 8928902
$ws.Range("J113").Value = 1338.3334
$ws.Range("K113").Value = 26786706
$ws.Range("L113").Value = 4015.0002
$ws.Range("M113").Value = -26784536
$ws.Range("N113").Value = -8355.0002
$ws.Range("H121").Value = 21858.928
$ws.Range("I121").Value = 559.25
$ws.Range("K121").Value = 1677.75
$ws.Range("M121").Value = -367.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124").Value = 34999
$ws.Range("J124").Value = 34999
$ws.Range("L124").Value = 34999
$ws.Range("N124").Value = -44819
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 5273.7
$ws.Range("I132").Value = 3580.9167
$ws.Range("K132").Value = 10742.7501
$ws.Range("M132").Value = -8212.750100000001
$ws.Range("H136").Value = 7363.391
$ws.Range("I136").Value = 6631.8335
$ws.Range("J136").Value = 9997
$ws.Range("K136").Value = 19895.5005
$ws.Range("L136").Value = 29991
$ws.Range("M136").Value = -17345.5005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 49999
$ws.Range("I49").Value = 49999
$ws.Range("K49").Value = 49999
$ws.Range("M49").Value = -49769
$ws.Range("H107").Value = 11765618
$ws.Range("I107").Value = 16667444
$ws.Range("J107").Value = 1235.6
$ws.Range("K107").Value = 50002332
$ws.Range("L107").Value = 3706.8
$ws.Range("M107").Value = -50000412
$ws.Range("N107").Value = -7546.799999999999
$ws.Range("H122").Value = 11000.904
$ws.Range("J122").Value = 19935.4
$ws.Range("L122").Value = 59806.2
$ws.Range("N122").Value = -64706.2
$ws.Range("H132").Value = 5232.0444
$ws.Range("I132").Value = 3534.658
$ws.Range("J132").Value = 14446.429
$ws.Range("K132").Value = 10603.974
$ws.Range("L132").Value = 43339.287
$ws.Range("M132").Value = -8073.974
$ws.Range("N132").Value = -48399.287
$ws.Range("H136").Value = 8779778
$ws.Range("I136").Value = 13165780
$ws.Range("K136").Value = 39497340
$ws.Range("M136").Value = -39494790
